$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("B5").AddComment("Test comment")
Write-Output "added"
Write-Output $ws.Range("B5").Comment.Text()
